# Statusbericht: "Naechste Schritte" (B19) und "Erledigte Meilensteine" (B21/B18)
# eingetragen -Anna

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Naechste Schritte -> Setup/Kickoff-Meilensteine (row 18/19 block)
$ws.Range("B19").Value = "Projektkoordination dokumentiert und abgeschlossen`nProjektcontrolling`nNetworking auf der Party betrieben`nProjekt erfolgreich abgeschlossen"

$ws.Range("B18").Value = "Design komplett und zur Einbindung in das Spiel übergeben`nPräsentables Spiel erzeugt`nFortschritt durch Projektleiter als zumindest ausreichend bewertet`nPräsentation präsentierfertig vorbereitet`nPräsentation erfolgreich absolviert"

# Erledigte Meilensteine
$ws.Range("B21").Value = "Grundidee gefunden `nProjekt Pitch Kick Off`nRequirements festgelegt und priorisiert`nTechstack spezifizieren/ Technologien ausgewählt`nSpezifikation von Projektleiter approved`nSetup abgeschlossen "

$ws.Range("A21").Select()
